$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("华银电力", "利欧股份", "豫能控股"),
    @("豫能控股", "华银电力", "寒武纪"),
    @("中天科技", "拓维信息", "利欧股份"),
    @("拓维信息", "包钢股份", "岩山科技"),
    @("章源钨业", "章源钨业", "拓维信息"),
    @("利欧股份", "豫能控股", "金正大"),
    @("包钢股份", "湖南黄金", "包钢股份"),
    @("金正大", "云南锗业", "华银电力"),
    @("岩山科技", "中天科技", "章源钨业"),
    @("湖南黄金", "金正大", "平潭发展"),
    @("平潭发展", "厦门钨业", "烽火通信"),
    @("烽火通信", "岩山科技", "万向钱潮"),
    @("协鑫集成", "中钨高新", "湖南黄金"),
    @("中钨高新", "北方稀土", "中钨高新"),
    @("首都在线", "平潭发展", "杭电股份"),
    @("寒武纪-U", "赣能股份", "协鑫集成"),
    @("云南锗业", "协鑫集成", "南网能源"),
    @("北方稀土", "航天发展", "大位科技"),
    @("航天发展", "烽火通信", "江钨装备")
)

$startRow = 3
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
